# Adding age at start of exp into excel dissec file
# Adds two new columns to the right of the existing data table:
#   I: "Pedigree prefix"  (mirrors the mouse ID already present in column D)
#   J: "Age start (d)"    (age, in days, of the animal at the start of the experiment)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "Pedigree prefix"
$ws.Range("J1").Value = "Age start (d)"

# row -> (pedigree prefix / ID, age start in days)
$data = @{
  2  = @("15833A", 83)
  3  = @("15894D", 75)
  5  = @("19017A", 64)
  6  = @("15842C", 85)
  7  = @("15978A", 67)
  8  = @("15935C", 68)
  9  = @("15913C", 77)
  10 = @("15820B", 88)
  11 = @("15834B", 83)
  12 = @("19005D", 64)
  13 = @("15980C", 67)
  14 = @("19018B", 64)
  15 = @("15891A", 75)
  16 = @("15822D", 88)
  18 = @("15970C", 72)
  19 = @("15841B", 85)
  20 = @("15899A", 73)
  21 = @("15845B", 85)
  22 = @("15947B", 65)
  23 = @("15857B", 81)
  24 = @("19004C", 64)
  25 = @("15905A", 78)
  26 = @("15956A", 70)
  27 = @("15840A", 85)
  30 = @("15844A", 85)
  31 = @("15999A", 64)
  32 = @("15862C", 81)
  33 = @("15948C", 65)
  34 = @("15856A", 81)
  35 = @("15843D", 85)
  36 = @("15983B", 67)
  37 = @("15912B", 77)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  $ws.Cells.Item($r, 9).Value = $vals[0]
  $ws.Cells.Item($r, 10).Value = $vals[1]
}

# Leave the selection where Excel would land after typing the last entry
[void]$ws.Range("I38").Select()
